$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 5.120597
$ws.Range("H2").Value = 15.361791
$ws.Range("I2").Value = 0.6410827008352843
$ws.Range("J2").Value = 0.6410827008352843
$ws.Range("M2").Value = 36.923013
$ws.Range("N2").Value = 110.769039
$ws.Range("O2").Value = 0.7437819354528793
$ws.Range("P2").Value = 0.7437819354528794
$ws.Range("Q2").Value = 189.067869598761
$ws.Range("R2").Value = 1701.610826388849
$ws.Range("S2").Value = 0.476825732012627
$ws.Range("T2").Value = 0.476825732012627
$ws.Range("G3").Value = 5.120597
$ws.Range("H3").Value = 15.361791
$ws.Range("I3").Value = 0.6410827008352843
$ws.Range("J3").Value = 0.6410827008352843
$ws.Range("O3").Value = 0.1364233939221953
$ws.Range("P3").Value = 0.1364233939221953
$ws.Range("Q3").Value = 34.678551901905
$ws.Range("R3").Value = 312.106967117145
$ws.Range("S3").Value = 0.08745867783275689
$ws.Range("T3").Value = 0.08745867783275689
$ws.Range("G4").Value = 5.120597
$ws.Range("H4").Value = 15.361791
$ws.Range("I4").Value = 0.6410827008352843
$ws.Range("J4").Value = 0.6410827008352843
$ws.Range("O4").Value = 0.1197946706249253
$ws.Range("P4").Value = 0.1197946706249254
$ws.Range("Q4").Value = 30.45156393930033
$ws.Range("R4").Value = 274.064075453703
$ws.Range("S4").Value = 0.07679829098990043
$ws.Range("T4").Value = 0.07679829098990044
$ws.Range("I5").Value = 0.3099611509948915
$ws.Range("J5").Value = 0.3099611509948915
$ws.Range("M5").Value = 36.923013
$ws.Range("N5").Value = 110.769039
$ws.Range("O5").Value = 0.7437819354528793
$ws.Range("P5").Value = 0.7437819354528794
$ws.Range("Q5").Value = 91.41362635527
$ws.Range("R5").Value = 822.7226371974299
$ws.Range("S5").Value = 0.2305435048021825
$ws.Range("T5").Value = 0.2305435048021826
$ws.Range("I6").Value = 0.3099611509948915
$ws.Range("J6").Value = 0.3099611509948915
$ws.Range("O6").Value = 0.1364233939221953
$ws.Range("P6").Value = 0.1364233939221953
$ws.Range("S6").Value = 0.04228595220275314
$ws.Range("T6").Value = 0.04228595220275314
$ws.Range("I7").Value = 0.3099611509948915
$ws.Range("J7").Value = 0.3099611509948915
$ws.Range("O7").Value = 0.1197946706249253
$ws.Range("P7").Value = 0.1197946706249254
$ws.Range("S7").Value = 0.03713169398995577
$ws.Range("T7").Value = 0.03713169398995578
$ws.Range("I8").Value = 0.04895614816982421
$ws.Range("J8").Value = 0.04895614816982421
$ws.Range("M8").Value = 36.923013
$ws.Range("N8").Value = 110.769039
$ws.Range("O8").Value = 0.7437819354528793
$ws.Range("P8").Value = 0.7437819354528794
$ws.Range("Q8").Value = 14.4381288501
$ws.Range("R8").Value = 129.9431596509
$ws.Range("S8").Value = 0.03641269863806978
$ws.Range("T8").Value = 0.03641269863806979
$ws.Range("I9").Value = 0.04895614816982421
$ws.Range("J9").Value = 0.04895614816982421
$ws.Range("O9").Value = 0.1364233939221953
$ws.Range("P9").Value = 0.1364233939221953
$ws.Range("R9").Value = 23.8339841445
$ws.Range("S9").Value = 0.006678763886685289
$ws.Range("T9").Value = 0.006678763886685289
$ws.Range("I10").Value = 0.04895614816982421
$ws.Range("J10").Value = 0.04895614816982421
$ws.Range("O10").Value = 0.1197946706249253
$ws.Range("P10").Value = 0.1197946706249254
$ws.Range("S10").Value = 0.005864685645069132
$ws.Range("T10").Value = 0.005864685645069133